# Regional statistics pipeline finished:
#  - Rename Sheet1 -> "country"
#  - Add a new sheet "region" (after "country") with code/alias data
#  - Make "region" the active sheet/tab
#  - Update the saved selections on both sheets

$wb = $excel.ActiveWorkbook

# Rename the original sheet.
$wsCountry = $wb.Worksheets.Item(1)
$wsCountry.Name = "country"

# Add the new "region" sheet right after "country".
$wsRegion = $wb.Worksheets.Add($null, $wsCountry)
$wsRegion.Name = "region"

# Populate the region sheet with its code/alias rows.
$wsRegion.Range("A1").Value = "code"
$wsRegion.Range("B1").Value = "alias"
$wsRegion.Range("A2").Value = "hlth_rs_physreg"
$wsRegion.Range("B2").Value = "physicians"
$wsRegion.Range("A3").Value = "tran_r_acci"
$wsRegion.Range("B3").Value = "road_accidents"

# Restore the last-used selections on each sheet.
$wsCountry.Range("C22").Select() | Out-Null
$wsRegion.Range("C18").Select() | Out-Null

# "region" is the active/visible tab when the workbook is saved.
$wsRegion.Activate() | Out-Null
